$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-36
# from 2026-02-22 (46075) to 2026-02-23 (46076)
$ws.Range("C2:C36").Value = 46076
